# "stop the train correctly red line" (#320)
# Move the Beacon (E) / B0 (F) / B1 (G) trigger values that were attached
# to the wrong block row up to the row that actually corresponds to the
# station named in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Red Line")

# --- Row 9 -> Row 8 (SHADYSIDE beacon) ---
$ws.Range("E8").Value2 = $ws.Range("E9").Value2
$ws.Range("F8").Value2 = $ws.Range("F9").Value2
$ws.Range("E9").Clear()
$ws.Range("F9").Clear()

# --- Row 23 -> Row 22 (SWISSVILLE beacon) ---
$ws.Range("E22").Value2 = $ws.Range("E23").Value2
$ws.Range("G22").Value2 = $ws.Range("G23").Value2
$ws.Range("E23").Clear()
$ws.Range("G23").Clear()

# --- Row 47 -> Row 46 (FIRST AVE beacon) ---
$ws.Range("E46").Value2 = $ws.Range("E47").Value2
$ws.Range("G46").Value2 = $ws.Range("G47").Value2
$ws.Range("E47").Clear()
$ws.Range("G47").Clear()

# --- Row 48 -> Row 49 (STATION SQUARE beacon) ---
$ws.Range("E49").Value2 = $ws.Range("E48").Value2
$ws.Range("F49").Value2 = $ws.Range("F48").Value2
$ws.Range("E48").Clear()
$ws.Range("F48").Clear()

# --- View state: restore to top of sheet with rows 12-17 frozen/visible ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select()
$ws.Range("A13").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A13").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E17").Select()
